$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Address" datatype block (row 7): Address / adr
$ws.Range("B7").Value = "Address"
$ws.Range("C7").Value = "adr"

# New "Environment" import block (rows 10-11)
$ws.Range("B10").Value = "Environment"
$ws.Range("B11").Value = "import"
$ws.Range("C11").Value = "com.example.beans"

# C10 stays empty but keeps the same formatted look as the rest of the
# block, so give it an explicit (touched) format like its neighbours.
$ws.Range("C10").Font.Bold = $false

# Column widths - best fit, matching content width
$ws.Columns("B").ColumnWidth = 14.75
$ws.Columns("C").ColumnWidth = 18.25
$ws.Columns("D").ColumnWidth = 24.09

# Selection as left by the editor
$ws.Range("C10").Select()
